$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update status of "Job list page needs customer details box" (row 21) to completed
$ws.Range("C21").Value = "Completed - March 29, 2010"

# Update the selected cell to C22 (as reflected in the saved view state)
$ws.Range("C22").Select()
